$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6337.074482800014
$ws.Range("D2").Value = 427.05887645
$ws.Range("B3").Value = 5915.804813150013
$ws.Range("D3").Value = 397.1648567
$ws.Range("B4").Value = 6308.487047833346
$ws.Range("D4").Value = 429.5158940333333
$ws.Range("B5").Value = 6109.31012581668
$ws.Range("D5").Value = 423.3215779
$ws.Range("B6").Value = 6333.656942483347
$ws.Range("D6").Value = 427.5835390833333
$ws.Range("B7").Value = 6112.18355251668
$ws.Range("D7").Value = 405.6004834166666
$ws.Range("B8").Value = 6302.481822583347
$ws.Range("D8").Value = 419.755824
$ws.Range("B9").Value = 6318.33731371668
$ws.Range("D9").Value = 444.3989366666667
$ws.Range("B10").Value = 6114.464430750013
$ws.Range("D10").Value = 415.7302043
$ws.Range("B11").Value = 6338.98541771668
$ws.Range("D11").Value = 421.5538573666667
$ws.Range("B12").Value = 6116.052410050012
$ws.Range("D12").Value = 416.1368788166666
$ws.Range("B13").Value = 6128.743814600013
$ws.Range("D13").Value = 418.08523845
